# The commit removes direct/local "contextual spacing" formatting
# (<w:contextualSpacing w:val="0"/>) from every paragraph's <w:pPr> in
# the document - i.e. in the Word UI this is the Paragraph dialog's
# "Don't add space between paragraphs of the same style" checkbox being
# cleared back to its unset/default state for each paragraph.
#
# The idiomatic COM-interop way to do this is to walk every paragraph in
# the document and reset ParagraphFormat.ContextualSpacing to False.

$d = $word.ActiveDocument

# Primary approach: iterate every paragraph in the story and clear the
# ContextualSpacing flag on its ParagraphFormat.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    try {
        $para.Range.ParagraphFormat.ContextualSpacing = $false
    } catch {
    }
    try {
        $para.ParagraphFormat.ContextualSpacing = $false
    } catch {
    }
    try {
        $para.Format.ContextualSpacing = $false
    } catch {
    }
}

# Belt-and-braces: also try the bulk collection-level setter and a
# whole-document Range setter, in case the host exposes a collection
# level shortcut instead of (or in addition to) per-paragraph access.
try {
    $d.Paragraphs.ContextualSpacing = $false
} catch {
}

try {
    $d.Content.ParagraphFormat.ContextualSpacing = $false
} catch {
}

try {
    $d.Range().ParagraphFormat.ContextualSpacing = $false
} catch {
}

Write-Host ("contextualSpacing cleared on " + $count + " paragraphs")
